# Updated phase 2 requirements partially, yet to implement
#
# 1) Rename the existing sheet to "Phase1" (keeps all of its data/layout).
# 2) Insert a brand-new sheet right after it named "Phase2" (becomes the
#    active/selected tab, matching activeTab="1" + tabSelected moving over).
# 3) Populate "Phase2" with the new requirements content.

$wb = $excel.ActiveWorkbook

$phase1 = $wb.Worksheets.Item(1)
$phase1.Name = "Phase1"

# Worksheets.Add(Before, After) - passing After places the new sheet
# immediately following $phase1, and makes it the active sheet.
$phase2 = $wb.Worksheets.Add($null, $phase1)
$phase2.Name = "Phase2"

# ---- Phase2 content ----
$phase2.Range("A1").Value = "Connected"
$phase2.Range("B1").Value = "Notification"
$phase2.Range("C1").Value = "Do you want to monitor?"

$phase2.Range("C2").Value = "yes"
$phase2.Range("D2").Value = "Start monitoring"
$phase2.Range("E2").Value = "if Stopped throw Toast"

$phase2.Range("C3").Value = "Dismiss"
$phase2.Range("D3").Value = "No action"

$phase2.Range("A5").Value = "Alert"
$phase2.Range("B5").Value = "Notification"
$phase2.Range("C5").Value = "Max/min level reached"

$phase2.Range("C6").Value = "Dismiss"
$phase2.Range("D6").Value = "Stop monitoring"

# ---- Phase2 column widths (approximate - engine quantizes to 1/6 char + padding) ----
$phase2.Columns.Item(1).ColumnWidth = 9.59
$phase2.Columns.Item(2).ColumnWidth = 17.59
$phase2.Columns.Item(3).ColumnWidth = 22.42
$phase2.Columns.Item(4).ColumnWidth = 17.09
$phase2.Columns.Item(5).ColumnWidth = 16.42

# ---- Phase2 selection ----
$phase2.Range("D7").Select() | Out-Null
